$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The instruction address/offset fields for the jump-class instructions
# (JUMP, JE, JA, JB, JBE, JAE — rows 10-15) are extended from 8 bits
# (columns J:Q) to the full 12 bits (columns F:Q), absorbing the 4-bit
# field that used to sit in columns F:I for those rows.
for ($r = 10; $r -le 15; $r++) {
    $label = $ws.Range("J$r").Value2

    # Grab the formatting of the existing label merge (bold font, outer
    # medium border, centered) before we touch anything.
    $ws.Range("J$r" + ":Q$r").Copy()
    $ws.Range("J$r" + ":Q$r").UnMerge()

    # Stamp that formatting across the whole new F:Q span.
    $ws.Range("F$r" + ":Q$r").PasteSpecial(-4122)
    $ws.Application.CutCopyMode = $false

    # Put the label back in F and clear the rest, then merge F:Q.
    $ws.Range("G$r" + ":Q$r").ClearContents()
    $ws.Range("F$r").Value2 = $label
    $ws.Range("F$r" + ":Q$r").Merge()
}

# Restore the selection the author left the sheet in.
$ws.Range("U6").Select()
